$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: LawDate (C) and LawDateMonth (D)
$ws.Range("C26").Value = 2019
$ws.Range("D26").Value = 201810

# Row 42: ExistingPolicy (B), LawDate (C), LawDateMonth (D)
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 2018
$ws.Range("D42").Value = 201805

# Row 51: LawDate (C), LawDateMonth (D)
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 0

# Update selection to C3
$ws.Activate()
$ws.Range("C3").Select()
